$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The classroom seating-plan generator was re-run ("Fixed Application Running
# Info"), producing a new arrangement of entities (pupils) across the two
# seating-plan sheets. Row 1 (titles) and the "zones" sheet are untouched;
# only the seat grids (rows 2-4, columns A-J) on seating_plan_maths and
# seating_plan_english change. The newly-populated seat A (front-left, ids
# "A3/A4") and the two spare seats at the end of the back row ("I4/J4") are
# now present but hold no entity, so they are written as blank placeholders.
# ---------------------------------------------------------------------------

$wsMaths = $wb.Worksheets.Item("seating_plan_maths")

$wsMaths.Range("A2").Value = "Katrina Petersone"
$wsMaths.Range("B2").Value = "Violet Hudson"
$wsMaths.Range("C2").Value = "Benjamin Finn"
$wsMaths.Range("D2").Value = "Samuel Dixon"
$wsMaths.Range("E2").Value = "Ava Lee"
$wsMaths.Range("F2").Value = "William Hunt"
$wsMaths.Range("G2").Value = "Caitlin Boyd"
$wsMaths.Range("H2").Value = "Madison Taylor"
$wsMaths.Range("I2").Value = "James Eilbeck"
$wsMaths.Range("J2").Value = "Ruby Haigh"

$wsMaths.Range("A3").Value = " "
$wsMaths.Range("B3").Value = "Nancy Enyoazu"
$wsMaths.Range("C3").Value = "James Calderon"
$wsMaths.Range("D3").Value = "Matthew Homan"
$wsMaths.Range("E3").Value = "Benedict Hobday"
$wsMaths.Range("F3").Value = "Alex Sentance"
$wsMaths.Range("G3").Value = "Thomas Barrett"
$wsMaths.Range("H3").Value = "James Shilton"
$wsMaths.Range("I3").Value = "Niko Morris"
$wsMaths.Range("J3").Value = "Lexi Green"

$wsMaths.Range("A4").Value = " "
$wsMaths.Range("B4").Value = "Aarron Kelly"
$wsMaths.Range("C4").Value = "Esther Sido"
$wsMaths.Range("D4").Value = "Brooke Layton"
$wsMaths.Range("E4").Value = "Sophie Rayner"
$wsMaths.Range("F4").Value = "Stanley Hirst"
$wsMaths.Range("G4").Value = "Benjamin Hillary"
$wsMaths.Range("H4").Value = "Spencer Rowe"
$wsMaths.Range("I4").Value = " "
$wsMaths.Range("J4").Value = " "

$wsEnglish = $wb.Worksheets.Item("seating_plan_english")

$wsEnglish.Range("A2").Value = "Callum Foster"
$wsEnglish.Range("B2").Value = "Cassie Strachan"
$wsEnglish.Range("C2").Value = "Ava Lee"
$wsEnglish.Range("D2").Value = "Jayden Nasa-Mereni"
$wsEnglish.Range("E2").Value = "Lewis Dacre"
$wsEnglish.Range("F2").Value = "Jayden Parsons"
$wsEnglish.Range("G2").Value = "Isabella Holmes"
$wsEnglish.Range("H2").Value = "Aarron Kelly"
$wsEnglish.Range("I2").Value = "Elliott Long"
$wsEnglish.Range("J2").Value = "Ethan Durham"

$wsEnglish.Range("A3").Value = " "
$wsEnglish.Range("B3").Value = "Lexie Starkey"
$wsEnglish.Range("C3").Value = "Caitlin Boyd"
$wsEnglish.Range("D3").Value = "Mariam Keita"
$wsEnglish.Range("E3").Value = "Niamh Teale"
$wsEnglish.Range("F3").Value = "James Eilbeck"
$wsEnglish.Range("G3").Value = "Bethany Greer"
$wsEnglish.Range("H3").Value = "James Calderon"
$wsEnglish.Range("I3").Value = "Samuel Dixon"
$wsEnglish.Range("J3").Value = "Cheryl Kanyimo"

$wsEnglish.Range("A4").Value = " "
$wsEnglish.Range("B4").Value = "Matthew Homan"
$wsEnglish.Range("C4").Value = "Alex Sentance"
$wsEnglish.Range("D4").Value = "Hugo Bird"
$wsEnglish.Range("E4").Value = "Eva Redican"
$wsEnglish.Range("F4").Value = "Lucy Webster"
$wsEnglish.Range("G4").Value = "Jude Fitzsimons"
$wsEnglish.Range("H4").Value = "Patryk Rudnicki"
$wsEnglish.Range("I4").Value = " "
$wsEnglish.Range("J4").Value = " "
